$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns at R:T, shifting existing R:T..AE columns right to U:AH.
# (xlShiftToRight = -4161)
$ws.Range("R1:T2").Insert(-4161)

# New header cells created by the insert (row 1)
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data cells created by the insert (row 2)
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Lowercase the "Unknown" placeholders for the generated columns D:J (state.State in K stays "Unknown")
$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "unknown"
$ws.Range("F2").Value = "unknown"
$ws.Range("G2").Value = "unknown"
$ws.Range("H2").Value = "unknown"
$ws.Range("I2").Value = "unknown"
$ws.Range("J2").Value = "unknown"

$wb.Save()
